# "Generate Report for Handback"
# A new handback run refreshed the timestamps recorded for the
# 59b68d80-0cb3-4ce2-9165-59d760900154 file across the Overview sheet and
# the per-locale (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-09-06 11:04:27"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-09-06 11:04:23"
$zhcn.Range("K2").Value = "2016-09-06 11:04:40"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-09-06 11:04:27"
$dede.Range("K2").Value = "2016-09-06 11:04:48"
